$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9655848145484924
$ws.Range("B1").Value = 2.123378753662109
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.823230743408203
$ws.Range("E1").Value = 1.118522882461548
